$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(42, 4).Value = 44231
$ws.Cells.Item(42, 10).Value = 70
$ws.Cells.Item(42, 11).Value = 7500
$ws.Cells.Item(42, 12).Value = 8000
$ws.Cells.Item(42, 13).Value = 7714
$ws.Cells.Item(42, 15).Value = "Región Metropolitana"
$ws.Cells.Item(42, 16).Value = 771

$ws.Cells.Item(43, 4).Value = 44421
$ws.Cells.Item(43, 10).Value = 100
$ws.Cells.Item(43, 11).Value = 7000
$ws.Cells.Item(43, 12).Value = 7500
$ws.Cells.Item(43, 13).Value = 7250
$ws.Cells.Item(43, 15).Value = "Región Metropolitana"
$ws.Cells.Item(43, 16).Value = 725

$ws.Cells.Item(44, 4).Value = 44336
$ws.Cells.Item(44, 10).Value = 60
$ws.Cells.Item(44, 11).Value = 6000
$ws.Cells.Item(44, 12).Value = 6500
$ws.Cells.Item(44, 13).Value = 6250
$ws.Cells.Item(44, 15).Value = "Región Metropolitana"
$ws.Cells.Item(44, 16).Value = 625

$ws.Cells.Item(45, 4).Value = 44665
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 6500
$ws.Cells.Item(45, 12).Value = 7000
$ws.Cells.Item(45, 13).Value = 6750
$ws.Cells.Item(45, 15).Value = "Región Metropolitana"
$ws.Cells.Item(45, 16).Value = 675

$ws.Cells.Item(46, 4).Value = 44246
$ws.Cells.Item(46, 10).Value = 60
$ws.Cells.Item(46, 11).Value = 9000
$ws.Cells.Item(46, 12).Value = 10000
$ws.Cells.Item(46, 13).Value = 9500
$ws.Cells.Item(46, 15).Value = "Región Metropolitana"
$ws.Cells.Item(46, 16).Value = 950

$ws.Cells.Item(47, 4).Value = 44469
$ws.Cells.Item(47, 10).Value = 100
$ws.Cells.Item(47, 11).Value = 6000
$ws.Cells.Item(47, 12).Value = 6500
$ws.Cells.Item(47, 13).Value = 6250
$ws.Cells.Item(47, 15).Value = "Región Metropolitana"
$ws.Cells.Item(47, 16).Value = 625

$ws.Cells.Item(48, 4).Value = 44741
$ws.Cells.Item(48, 10).Value = 100
$ws.Cells.Item(48, 11).Value = 8000
$ws.Cells.Item(48, 12).Value = 8500
$ws.Cells.Item(48, 13).Value = 8250
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 825

$ws.Cells.Item(49, 4).Value = 44376
$ws.Cells.Item(49, 10).Value = 100
$ws.Cells.Item(49, 11).Value = 6000
$ws.Cells.Item(49, 12).Value = 6500
$ws.Cells.Item(49, 13).Value = 6250
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 625

$ws.Cells.Item(50, 4).Value = 44755
$ws.Cells.Item(50, 10).Value = 40
$ws.Cells.Item(50, 11).Value = 9000
$ws.Cells.Item(50, 12).Value = 10000
$ws.Cells.Item(50, 13).Value = 9500
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 950

$ws.Cells.Item(51, 4).Value = 44580
$ws.Cells.Item(51, 10).Value = 100
$ws.Cells.Item(51, 11).Value = 10000
$ws.Cells.Item(51, 12).Value = 11000
$ws.Cells.Item(51, 13).Value = 10500
$ws.Cells.Item(51, 15).Value = "Región Metropolitana"
$ws.Cells.Item(51, 16).Value = 1050

$ws.Cells.Item(52, 4).Value = 44635
$ws.Cells.Item(52, 10).Value = 170
$ws.Cells.Item(52, 11).Value = 7500
$ws.Cells.Item(52, 12).Value = 8000
$ws.Cells.Item(52, 13).Value = 7765
$ws.Cells.Item(52, 15).Value = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value = 776

$ws.Cells.Item(53, 4).Value = 44776
$ws.Cells.Item(53, 10).Value = 150
$ws.Cells.Item(53, 11).Value = 6500
$ws.Cells.Item(53, 12).Value = 7000
$ws.Cells.Item(53, 13).Value = 6733
$ws.Cells.Item(53, 15).Value = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value = 673

$ws.Cells.Item(54, 4).Value = 44313
$ws.Cells.Item(54, 10).Value = 60
$ws.Cells.Item(54, 11).Value = 6000
$ws.Cells.Item(54, 12).Value = 6500
$ws.Cells.Item(54, 13).Value = 6250
$ws.Cells.Item(54, 15).Value = "Región Metropolitana"
$ws.Cells.Item(54, 16).Value = 625

$ws.Cells.Item(55, 4).Value = 44238
$ws.Cells.Item(55, 10).Value = 100
$ws.Cells.Item(55, 11).Value = 8000
$ws.Cells.Item(55, 12).Value = 8500
$ws.Cells.Item(55, 13).Value = 8250
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 825

$ws.Cells.Item(56, 4).Value = 44672
$ws.Cells.Item(56, 10).Value = 220
$ws.Cells.Item(56, 11).Value = 10000
$ws.Cells.Item(56, 12).Value = 11000
$ws.Cells.Item(56, 13).Value = 10455
$ws.Cells.Item(56, 15).Value = "Región Metropolitana"
$ws.Cells.Item(56, 16).Value = 1046

$ws.Cells.Item(57, 4).Value = 44355
$ws.Cells.Item(57, 10).Value = 50
$ws.Cells.Item(57, 11).Value = 6000
$ws.Cells.Item(57, 12).Value = 6500
$ws.Cells.Item(57, 13).Value = 6300
$ws.Cells.Item(57, 15).Value = "Región Metropolitana"
$ws.Cells.Item(57, 16).Value = 630

$ws.Cells.Item(58, 4).Value = 44299
$ws.Cells.Item(58, 10).Value = 100
$ws.Cells.Item(58, 11).Value = 8000
$ws.Cells.Item(58, 12).Value = 9000
$ws.Cells.Item(58, 13).Value = 8500
$ws.Cells.Item(58, 15).Value = "Región Metropolitana"
$ws.Cells.Item(58, 16).Value = 850

$ws.Cells.Item(59, 4).Value = 44292
$ws.Cells.Item(59, 10).Value = 50
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 11000
$ws.Cells.Item(59, 13).Value = 10600
$ws.Cells.Item(59, 15).Value = "Región Metropolitana"
$ws.Cells.Item(59, 16).Value = 1060

$ws.Cells.Item(60, 4).Value = 44308
$ws.Cells.Item(60, 10).Value = 100
$ws.Cells.Item(60, 11).Value = 5000
$ws.Cells.Item(60, 12).Value = 5500
$ws.Cells.Item(60, 13).Value = 5250
$ws.Cells.Item(60, 15).Value = "Región Metropolitana"
$ws.Cells.Item(60, 16).Value = 525

$ws.Cells.Item(61, 4).Value = 44747
$ws.Cells.Item(61, 10).Value = 150
$ws.Cells.Item(61, 11).Value = 8000
$ws.Cells.Item(61, 12).Value = 8500
$ws.Cells.Item(61, 13).Value = 8233
$ws.Cells.Item(61, 15).Value = "Región Metropolitana"
$ws.Cells.Item(61, 16).Value = 823

$ws.Cells.Item(62, 4).Value = 44692
$ws.Cells.Item(62, 10).Value = 100
$ws.Cells.Item(62, 11).Value = 6000
$ws.Cells.Item(62, 12).Value = 6500
$ws.Cells.Item(62, 13).Value = 6250
$ws.Cells.Item(62, 15).Value = "Región Metropolitana"
$ws.Cells.Item(62, 16).Value = 625

$ws.Cells.Item(63, 4).Value = 44719
$ws.Cells.Item(63, 10).Value = 100
$ws.Cells.Item(63, 11).Value = 7000
$ws.Cells.Item(63, 12).Value = 7500
$ws.Cells.Item(63, 13).Value = 7250
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 725

$ws.Cells.Item(64, 4).Value = 44489
$ws.Cells.Item(64, 10).Value = 50
$ws.Cells.Item(64, 11).Value = 6000
$ws.Cells.Item(64, 12).Value = 6500
$ws.Cells.Item(64, 13).Value = 6300
$ws.Cells.Item(64, 15).Value = "Región Metropolitana"
$ws.Cells.Item(64, 16).Value = 630

$ws.Cells.Item(65, 4).Value = 44160
$ws.Cells.Item(65, 10).Value = 100
$ws.Cells.Item(65, 11).Value = 9000
$ws.Cells.Item(65, 12).Value = 9500
$ws.Cells.Item(65, 13).Value = 9250
$ws.Cells.Item(65, 15).Value = "Región Metropolitana"
$ws.Cells.Item(65, 16).Value = 925

$ws.Cells.Item(66, 4).Value = 44433
$ws.Cells.Item(66, 10).Value = 100
$ws.Cells.Item(66, 11).Value = 7000
$ws.Cells.Item(66, 12).Value = 7500
$ws.Cells.Item(66, 13).Value = 7250
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value = 725

$ws.Cells.Item(67, 4).Value = 44825
$ws.Cells.Item(67, 10).Value = 50
$ws.Cells.Item(67, 11).Value = 7500
$ws.Cells.Item(67, 12).Value = 8000
$ws.Cells.Item(67, 13).Value = 7700
$ws.Cells.Item(67, 15).Value = "Región Metropolitana"
$ws.Cells.Item(67, 16).Value = 770

$ws.Cells.Item(68, 4).Value = 44565
$ws.Cells.Item(68, 10).Value = 40
$ws.Cells.Item(68, 11).Value = 12000
$ws.Cells.Item(68, 12).Value = 13000
$ws.Cells.Item(68, 13).Value = 12500
$ws.Cells.Item(68, 15).Value = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value = 1250

$ws.Cells.Item(69, 4).Value = 44204
$ws.Cells.Item(69, 10).Value = 80
$ws.Cells.Item(69, 11).Value = 7000
$ws.Cells.Item(69, 12).Value = 7500
$ws.Cells.Item(69, 13).Value = 7188
$ws.Cells.Item(69, 15).Value = "Región Metropolitana"
$ws.Cells.Item(69, 16).Value = 719

$ws.Cells.Item(70, 4).Value = 44194
$ws.Cells.Item(70, 10).Value = 100
$ws.Cells.Item(70, 11).Value = 8000
$ws.Cells.Item(70, 12).Value = 9000
$ws.Cells.Item(70, 13).Value = 8500
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 850

$ws.Cells.Item(71, 4).Value = 44645
$ws.Cells.Item(71, 10).Value = 160
$ws.Cells.Item(71, 11).Value = 6000
$ws.Cells.Item(71, 12).Value = 6500
$ws.Cells.Item(71, 13).Value = 6250
$ws.Cells.Item(71, 15).Value = "Región Metropolitana"
$ws.Cells.Item(71, 16).Value = 625

$ws.Cells.Item(72, 4).Value = 44526
$ws.Cells.Item(72, 10).Value = 40
$ws.Cells.Item(72, 11).Value = 8000
$ws.Cells.Item(72, 12).Value = 8500
$ws.Cells.Item(72, 13).Value = 8250
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 825

$ws.Cells.Item(73, 4).Value = 44383
$ws.Cells.Item(73, 10).Value = 60
$ws.Cells.Item(73, 11).Value = 7500
$ws.Cells.Item(73, 12).Value = 8000
$ws.Cells.Item(73, 13).Value = 7750
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 775

$ws.Cells.Item(74, 4).Value = 44784
$ws.Cells.Item(74, 10).Value = 40
$ws.Cells.Item(74, 11).Value = 6500
$ws.Cells.Item(74, 12).Value = 7000
$ws.Cells.Item(74, 13).Value = 6750
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 675

$ws.Cells.Item(75, 4).Value = 44166
$ws.Cells.Item(75, 10).Value = 100
$ws.Cells.Item(75, 11).Value = 8000
$ws.Cells.Item(75, 12).Value = 9000
$ws.Cells.Item(75, 13).Value = 8500
$ws.Cells.Item(75, 15).Value = "Región Metropolitana"
$ws.Cells.Item(75, 16).Value = 850

$ws.Cells.Item(76, 4).Value = 44582
$ws.Cells.Item(76, 10).Value = 50
$ws.Cells.Item(76, 11).Value = 8000
$ws.Cells.Item(76, 12).Value = 8500
$ws.Cells.Item(76, 13).Value = 8200
$ws.Cells.Item(76, 15).Value = "Región Metropolitana"
$ws.Cells.Item(76, 16).Value = 820

$ws.Cells.Item(77, 4).Value = 44881
$ws.Cells.Item(77, 10).Value = 50
$ws.Cells.Item(77, 11).Value = 13000
$ws.Cells.Item(77, 12).Value = 14000
$ws.Cells.Item(77, 13).Value = 13400
$ws.Cells.Item(77, 15).Value = "Región Metropolitana"
$ws.Cells.Item(77, 16).Value = 1340

$ws.Cells.Item(78, 4).Value = 44476
$ws.Cells.Item(78, 10).Value = 80
$ws.Cells.Item(78, 11).Value = 5000
$ws.Cells.Item(78, 12).Value = 5500
$ws.Cells.Item(78, 13).Value = 5312
$ws.Cells.Item(78, 15).Value = "Región Metropolitana"
$ws.Cells.Item(78, 16).Value = 531

$ws.Cells.Item(79, 4).Value = 44923
$ws.Cells.Item(79, 10).Value = 100
$ws.Cells.Item(79, 11).Value = 7000
$ws.Cells.Item(79, 12).Value = 7500
$ws.Cells.Item(79, 13).Value = 7250
$ws.Cells.Item(79, 15).Value = "Región Metropolitana"
$ws.Cells.Item(79, 16).Value = 725

$ws.Cells.Item(80, 4).Value = 44923
$ws.Cells.Item(80, 10).Value = 100
$ws.Cells.Item(80, 11).Value = 7000
$ws.Cells.Item(80, 12).Value = 7500
$ws.Cells.Item(80, 13).Value = 7250
$ws.Cells.Item(80, 15).Value = "Región Metropolitana"
$ws.Cells.Item(80, 16).Value = 725

$ws.Cells.Item(81, 4).Value = 44209
$ws.Cells.Item(81, 10).Value = 80
$ws.Cells.Item(81, 11).Value = 7500
$ws.Cells.Item(81, 12).Value = 8000
$ws.Cells.Item(81, 13).Value = 7688
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 769

$ws.Cells.Item(82, 4).Value = 44264
$ws.Cells.Item(82, 10).Value = 50
$ws.Cells.Item(82, 11).Value = 8000
$ws.Cells.Item(82, 12).Value = 8500
$ws.Cells.Item(82, 13).Value = 8200
$ws.Cells.Item(82, 15).Value = "Región Metropolitana"
$ws.Cells.Item(82, 16).Value = 820

$ws.Cells.Item(83, 4).Value = 44782
$ws.Cells.Item(83, 10).Value = 100
$ws.Cells.Item(83, 11).Value = 7000
$ws.Cells.Item(83, 12).Value = 7500
$ws.Cells.Item(83, 13).Value = 7250
$ws.Cells.Item(83, 15).Value = "Región Metropolitana"
$ws.Cells.Item(83, 16).Value = 725

$ws.Cells.Item(84, 4).Value = 44951
$ws.Cells.Item(84, 10).Value = 50
$ws.Cells.Item(84, 11).Value = 7000
$ws.Cells.Item(84, 12).Value = 7500
$ws.Cells.Item(84, 13).Value = 7300
$ws.Cells.Item(84, 15).Value = "Región Metropolitana"
$ws.Cells.Item(84, 16).Value = 730
